$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1580.5
$ws.Range("I18").Value = 1440
$ws.Range("J18").Value = 2002
$ws.Range("K18").Value = 1440
$ws.Range("L18").Value = 2002
$ws.Range("M18").Value = -1156
$ws.Range("N18").Value = -2570
$ws.Range("H33").Value = 259.83334
$ws.Range("I33").Value = 290.8889
$ws.Range("K33").Value = 290.8889
$ws.Range("M33").Value = -61.88889999999998
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H74").Value = 12000
$ws.Range("I74").Value = 12000
$ws.Range("K74").Value = 12000
$ws.Range("M74").Value = -11064
$ws.Range("H77").Value = 12000
$ws.Range("I77").Value = 12000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55320
$ws.Range("H86").Value = 4840
$ws.Range("I86").Value = 3800
$ws.Range("K86").Value = 3800
$ws.Range("M86").Value = -2677
$ws.Range("H89").Value = 4840
$ws.Range("I89").Value = 3800
$ws.Range("K89").Value = 19000
$ws.Range("M89").Value = -13384
$ws.Range("H113").Value = 4399.6
$ws.Range("J113").Value = 4999.6665
$ws.Range("L113").Value = 4999.6665
$ws.Range("N113").Value = -11507.6665
$ws.Range("H125").Value = 905
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
$ws.Range("H137").Value = 2673.44
$ws.Range("I137").Value = 1695.3846
$ws.Range("J137").Value = 3733
$ws.Range("K137").Value = 5086.1538
$ws.Range("L137").Value = 11199
$ws.Range("M137").Value = -2536.1538
$ws.Range("N137").Value = -16299
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 83.8
$ws.Range("I5").Value = 23.333334
$ws.Range("J5").Value = 98.916664
$ws.Range("K5").Value = 23.333334
$ws.Range("L5").Value = 98.916664
$ws.Range("M5").Value = 88.66666599999999
$ws.Range("N5").Value = -322.916664
$ws.Range("H61").Value = 3784.875
$ws.Range("I61").Value = 3784.875
$ws.Range("K61").Value = 3784.875
$ws.Range("M61").Value = -3572.875
$ws.Range("H88").Value = 2318.5
$ws.Range("I88").Value = 2605
$ws.Range("J88").Value = 1841
$ws.Range("K88").Value = 2605
$ws.Range("L88").Value = 1841
$ws.Range("M88").Value = -2199
$ws.Range("N88").Value = -2653
$ws.Range("H91").Value = 2318.5
$ws.Range("I91").Value = 2605
$ws.Range("J91").Value = 1841
$ws.Range("K91").Value = 2605
$ws.Range("L91").Value = 1841
$ws.Range("M91").Value = -1201
$ws.Range("N91").Value = -4649
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 5607.615
$ws.Range("I132").Value = 5408.25
$ws.Range("K132").Value = 16224.75
$ws.Range("M132").Value = -13694.75
$ws.Range("H136").Value = 3784.875
$ws.Range("I136").Value = 3784.875
$ws.Range("K136").Value = 11354.625
$ws.Range("M136").Value = -8804.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 83.8
$ws.Range("I4").Value = 23.333334
$ws.Range("J4").Value = 98.916664
$ws.Range("K4").Value = 23.333334
$ws.Range("L4").Value = 98.916664
$ws.Range("M4").Value = 91.66666599999999
$ws.Range("N4").Value = -328.916664
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H102").Value = 9599.5
$ws.Range("I102").Value = 9599.5
$ws.Range("K102").Value = 9599.5
$ws.Range("M102").Value = -6354.5
$ws.Range("H134").Value = 1211
$ws.Range("I134").Value = 1211
$ws.Range("K134").Value = 3633
$ws.Range("M134").Value = -1098
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 11999.2
$ws.Range("J28").Value = 11999.2
$ws.Range("L28").Value = 11999.2
$ws.Range("N28").Value = -12489.2
$ws.Range("H86").Value = 4838.3335
$ws.Range("I86").Value = 2253.5
$ws.Range("J86").Value = 10008
$ws.Range("K86").Value = 2253.5
$ws.Range("L86").Value = 10008
$ws.Range("M86").Value = -1130.5
$ws.Range("N86").Value = -12254
$ws.Range("H89").Value = 4838.3335
$ws.Range("I89").Value = 2253.5
$ws.Range("J89").Value = 10008
$ws.Range("K89").Value = 11267.5
$ws.Range("L89").Value = 50040
$ws.Range("M89").Value = -5651.5
$ws.Range("N89").Value = -61272
$ws.Range("H94").Value = 4158.231
$ws.Range("I94").Value = 2253.5
$ws.Range("K94").Value = 2253.5
$ws.Range("M94").Value = -1802.5
$ws.Range("H122").Value = 1357.1428
$ws.Range("I122").Value = 1280.8
$ws.Range("K122").Value = 3842.4
$ws.Range("M122").Value = -1392.4
$ws.Range("H132").Value = 3679.2354
$ws.Range("I132").Value = 1273.6666
$ws.Range("J132").Value = 6385.5
$ws.Range("K132").Value = 3820.9998
$ws.Range("L132").Value = 19156.5
$ws.Range("M132").Value = -1290.9998
$ws.Range("N132").Value = -24216.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 87216.25999999999
$ws.Range("I4").Value = 142963.72
$ws.Range("K4").Value = 428891.16
$ws.Range("M4").Value = -428779.16
$ws.Range("H37").Value = 179999
$ws.Range("J37").Value = 179999
$ws.Range("L37").Value = 539997
$ws.Range("N37").Value = -540221
$ws.Range("H107").Value = 317.66666
$ws.Range("I107").Value = 351.5
$ws.Range("K107").Value = 1054.5
$ws.Range("M107").Value = 865.5
$ws.Range("H109").Value = 1566.3334
$ws.Range("I109").Value = 1633.3334
$ws.Range("J109").Value = 1499.3334
$ws.Range("K109").Value = 4900.0002
$ws.Range("L109").Value = 4498.0002
$ws.Range("M109").Value = -3860.0002
$ws.Range("N109").Value = -6578.0002
$ws.Range("H112").Value = 1227
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H121").Value = 336.75
$ws.Range("I121").Value = 282.33334
$ws.Range("J121").Value = 500
$ws.Range("K121").Value = 847.0000200000001
$ws.Range("L121").Value = 1500
$ws.Range("M121").Value = 462.9999799999999
$ws.Range("N121").Value = -4120
$ws.Range("H131").Value = 1813.0526
$ws.Range("J131").Value = 2243.3333
$ws.Range("L131").Value = 6729.999899999999
$ws.Range("N131").Value = -16809.9999
$ws.Range("H139").Value = 2371.7058
$ws.Range("I139").Value = 1515.6428
$ws.Range("K139").Value = 4546.928400000001
$ws.Range("M139").Value = 593.0715999999993
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6504
$ws.Range("I70").Value = 6504
$ws.Range("K70").Value = 6504
$ws.Range("M70").Value = -6234
$ws.Range("H73").Value = 6504
$ws.Range("I73").Value = 6504
$ws.Range("K73").Value = 6504
$ws.Range("M73").Value = -5568
$ws.Range("H80").Value = 3074.2222
$ws.Range("I80").Value = 2916.3333
$ws.Range("J80").Value = 3390
$ws.Range("K80").Value = 2916.3333
$ws.Range("L80").Value = 3390
$ws.Range("M80").Value = -1918.3333
$ws.Range("N80").Value = -5386
$ws.Range("H83").Value = 3074.2222
$ws.Range("I83").Value = 2916.3333
$ws.Range("J83").Value = 3390
$ws.Range("K83").Value = 14581.6665
$ws.Range("L83").Value = 16950
$ws.Range("M83").Value = -9589.666499999999
$ws.Range("N83").Value = -26934
$ws.Range("H132").Value = 6999.5
$ws.Range("J132").Value = 6999.5
$ws.Range("L132").Value = 20998.5
$ws.Range("N132").Value = -26058.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 10000
$ws.Range("K40").Value = 10000
$ws.Range("M40").Value = -9864
$ws.Range("H122").Value = 3348.6
$ws.Range("I122").Value = 3348.6
$ws.Range("K122").Value = 10045.8
$ws.Range("M122").Value = -7595.799999999999
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3360
$ws.Range("I132").Value = 3200
$ws.Range("K132").Value = 9600
$ws.Range("M132").Value = -7070
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 7500
$ws.Range("K136").Value = 22500
$ws.Range("M136").Value = -19950
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 12376
$ws.Range("I58").Value = 12376
$ws.Range("K58").Value = 12376
$ws.Range("M58").Value = -12068
$ws.Range("H81").Value = 5634.4546
$ws.Range("I81").Value = 2997.6667
$ws.Range("J81").Value = 17500
$ws.Range("K81").Value = 5995.3334
$ws.Range("L81").Value = 35000
$ws.Range("M81").Value = -4934.3334
$ws.Range("N81").Value = -37122
$ws.Range("H84").Value = 5634.4546
$ws.Range("I84").Value = 2997.6667
$ws.Range("J84").Value = 17500
$ws.Range("K84").Value = 29976.667
$ws.Range("L84").Value = 175000
$ws.Range("M84").Value = -24672.667
$ws.Range("N84").Value = -185608
$ws.Range("H132").Value = 2337.913
$ws.Range("I132").Value = 2363.7222
$ws.Range("K132").Value = 7091.1666
$ws.Range("M132").Value = -4561.1666
